$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6126.2144
$ws.Range("I32").Value = 6925.636
$ws.Range("J32").Value = 5608.9414
$ws.Range("K32").Value = 6925.636
$ws.Range("L32").Value = 5608.9414
$ws.Range("M32").Value = -6599.636
$ws.Range("N32").Value = -6260.9414

$ws.Range("H33").Value = 317.13043
$ws.Range("I33").Value = 215.05882
$ws.Range("J33").Value = 606.3333
$ws.Range("K33").Value = 215.05882
$ws.Range("L33").Value = 606.3333
$ws.Range("M33").Value = 13.94118
$ws.Range("N33").Value = -1064.3333

$ws.Range("H70").Value = 6948983.5
$ws.Range("I70").Value = 11112487
$ws.Range("J70").Value = 9810.333000000001
$ws.Range("K70").Value = 33337461
$ws.Range("L70").Value = 29430.999
$ws.Range("M70").Value = -33337191
$ws.Range("N70").Value = -29970.999

$ws.Range("H73").Value = 6948983.5
$ws.Range("I73").Value = 11112487
$ws.Range("J73").Value = 9810.333000000001
$ws.Range("K73").Value = 33337461
$ws.Range("L73").Value = 29430.999
$ws.Range("M73").Value = -33336525
$ws.Range("N73").Value = -31302.999

$ws.Range("H98").Value = 4999.75
$ws.Range("I98").Value = 2499.1
$ws.Range("J98").Value = 17503
$ws.Range("K98").Value = 2499.1
$ws.Range("L98").Value = 17503
$ws.Range("M98").Value = -1001.1
$ws.Range("N98").Value = -20499

$ws.Range("H113").Value = 8926.467000000001
$ws.Range("I113").Value = 8589.700000000001
$ws.Range("J113").Value = 9600
$ws.Range("K113").Value = 8589.700000000001
$ws.Range("L113").Value = 9600
$ws.Range("M113").Value = -5335.700000000001
$ws.Range("N113").Value = -16108

$ws.Range("H116").Value = 19227.809
$ws.Range("I116").Value = 21743.75
$ws.Range("K116").Value = 21743.75
$ws.Range("M116").Value = -18301.75

$ws.Range("H122").Value = 4999.75
$ws.Range("I122").Value = 2499.1
$ws.Range("J122").Value = 17503
$ws.Range("K122").Value = 7497.299999999999
$ws.Range("L122").Value = 52509
$ws.Range("M122").Value = -5047.299999999999
$ws.Range("N122").Value = -57409

$ws.Range("H132").Value = 25636.232
$ws.Range("I132").Value = 30945.75
$ws.Range("J132").Value = 4398.1665
$ws.Range("K132").Value = 92837.25
$ws.Range("L132").Value = 13194.4995
$ws.Range("M132").Value = -90307.25
$ws.Range("N132").Value = -18254.4995

$ws.Range("H137").Value = 26623.75
$ws.Range("I137").Value = 60799.4
$ws.Range("K137").Value = 182398.2
$ws.Range("M137").Value = -179848.2

$ws.Range("H138").Value = 2157.641
$ws.Range("I138").Value = 1868.8077
$ws.Range("J138").Value = 2735.3076
$ws.Range("K138").Value = 5606.4231
$ws.Range("L138").Value = 8205.9228
$ws.Range("M138").Value = -466.4231
$ws.Range("N138").Value = -18485.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19057.357
$ws.Range("I32").Value = 20051.207
$ws.Range("K32").Value = 20051.207
$ws.Range("M32").Value = -19764.207

$ws.Range("H61").Value = 8486.071
$ws.Range("I61").Value = 908.2727
$ws.Range("J61").Value = 36271.332
$ws.Range("K61").Value = 908.2727
$ws.Range("L61").Value = 36271.332
$ws.Range("M61").Value = -696.2727
$ws.Range("N61").Value = -36695.332

$ws.Range("H97").Value = 684.7059
$ws.Range("I97").Value = 664.931
$ws.Range("K97").Value = 664.931
$ws.Range("M97").Value = -168.931

$ws.Range("H102").Value = 3150.2222
$ws.Range("I102").Value = 3765.8572
$ws.Range("K102").Value = 3765.8572
$ws.Range("M102").Value = -2143.8572

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H132").Value = 2029.0344
$ws.Range("I132").Value = 1897.238
$ws.Range("K132").Value = 5691.714
$ws.Range("M132").Value = -3161.714

$ws.Range("H136").Value = 8486.071
$ws.Range("I136").Value = 908.2727
$ws.Range("J136").Value = 36271.332
$ws.Range("K136").Value = 2724.8181
$ws.Range("L136").Value = 108813.996
$ws.Range("M136").Value = -174.8181
$ws.Range("N136").Value = -113913.996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1104.1428
$ws.Range("I86").Value = 915.5454999999999
$ws.Range("J86").Value = 1311.6
$ws.Range("K86").Value = 915.5454999999999
$ws.Range("L86").Value = 1311.6
$ws.Range("M86").Value = 207.4545000000001
$ws.Range("N86").Value = -3557.6

$ws.Range("H89").Value = 1104.1428
$ws.Range("I89").Value = 915.5454999999999
$ws.Range("J89").Value = 1311.6
$ws.Range("K89").Value = 4577.7275
$ws.Range("L89").Value = 6558
$ws.Range("M89").Value = 1038.2725
$ws.Range("N89").Value = -17790

$ws.Range("H94").Value = 497.5
$ws.Range("I94").Value = 497.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 497.5
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -46.5
$ws.Range("N94").ClearContents()

$ws.Range("H99").Value = 1539.1
$ws.Range("I99").Value = 1298.875
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 1298.875
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = 199.125
$ws.Range("N99").Value = -5496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3572890
$ws.Range("J31").Value = 4399.8
$ws.Range("L31").Value = 4399.8
$ws.Range("N31").Value = -4989.8

$ws.Range("H34").Value = 3572890
$ws.Range("J34").Value = 4399.8
$ws.Range("L34").Value = 4399.8
$ws.Range("N34").Value = -4803.8

$ws.Range("H99").Value = 5220.7144
$ws.Range("J99").Value = 8850
$ws.Range("L99").Value = 8850
$ws.Range("N99").Value = -11846

$ws.Range("H126").Value = 5220.7144
$ws.Range("J126").Value = 8850
$ws.Range("L126").Value = 26550
$ws.Range("N126").Value = -31490

$ws.Range("H127").Value = 55000
$ws.Range("J127").Value = 55000
$ws.Range("L127").Value = 55000
$ws.Range("N127").Value = -64920

$ws.Range("H140").Value = 106653.47
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 106653.47
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 106653.47
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -117013.47

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 11344
$ws.Range("J19").Value = 11344
$ws.Range("L19").Value = 11344
$ws.Range("N19").Value = -11920

$ws.Range("H132").Value = 2922.6667
$ws.Range("I132").Value = 2562.6365
$ws.Range("K132").Value = 7687.9095
$ws.Range("M132").Value = -5157.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2832.65
$ws.Range("I55").Value = 1785.2
$ws.Range("J55").Value = 3880.1
$ws.Range("K55").Value = 1785.2
$ws.Range("L55").Value = 3880.1
$ws.Range("M55").Value = -1612.2
$ws.Range("N55").Value = -4226.1

$ws.Range("H61").Value = 3282.4546
$ws.Range("I61").Value = 3110.7
$ws.Range("K61").Value = 3110.7
$ws.Range("M61").Value = -2908.7

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H113").Value = 3282.4546
$ws.Range("I113").Value = 3110.7
$ws.Range("K113").Value = 3110.7
$ws.Range("M113").Value = -940.6999999999998

$ws.Range("H122").Value = 3886.6
$ws.Range("I122").Value = 3900
$ws.Range("K122").Value = 11700
$ws.Range("M122").Value = -9250

$ws.Range("H131").Value = 36485.75
$ws.Range("J131").Value = 31882.334
$ws.Range("L131").Value = 31882.334
$ws.Range("N131").Value = -41962.334

$ws.Range("H132").Value = 3975.5
$ws.Range("I132").Value = 3605.2727
$ws.Range("J132").Value = 5333
$ws.Range("K132").Value = 10815.8181
$ws.Range("L132").Value = 15999
$ws.Range("M132").Value = -8285.8181
$ws.Range("N132").Value = -21059

$ws.Range("H136").Value = 3509.8696
$ws.Range("I136").Value = 2588.8572
$ws.Range("J136").Value = 4942.5557
$ws.Range("K136").Value = 7766.571599999999
$ws.Range("L136").Value = 14827.6671
$ws.Range("M136").Value = -5216.571599999999
$ws.Range("N136").Value = -19927.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 8000000
$ws.Range("J6").Value = 6000000
$ws.Range("L6").Value = 6000000
$ws.Range("N6").Value = -6000230

$ws.Range("H100").Value = 593.6667
$ws.Range("I100").Value = 613
$ws.Range("J100").Value = 555
$ws.Range("K100").Value = 1226
$ws.Range("L100").Value = 1110
$ws.Range("M100").Value = -685
$ws.Range("N100").Value = -2192

$ws.Range("H120").Value = 60000
$ws.Range("J120").Value = 60000
$ws.Range("L120").Value = 60000
$ws.Range("N120").Value = -69676

$ws.Range("H124").Value = 57534.5
$ws.Range("J124").Value = 57534.5
$ws.Range("L124").Value = 57534.5
$ws.Range("N124").Value = -67354.5

$ws.Range("H126").Value = 253104.75
$ws.Range("I126").Value = 2462.1667
$ws.Range("J126").Value = 629068.6
$ws.Range("K126").Value = 7386.500100000001
$ws.Range("L126").Value = 1887205.8
$ws.Range("M126").Value = -4916.500100000001
$ws.Range("N126").Value = -1892145.8

$ws.Range("H132").Value = 1725.3889
$ws.Range("I132").Value = 1276.7142
$ws.Range("J132").Value = 3295.75
$ws.Range("K132").Value = 3830.1426
$ws.Range("L132").Value = 9887.25
$ws.Range("M132").Value = -1300.1426
$ws.Range("N132").Value = -14947.25

$ws.Range("H136").Value = 17702.086
$ws.Range("I136").Value = 26020.727
$ws.Range("J136").Value = 3624.3845
$ws.Range("K136").Value = 78062.181
$ws.Range("L136").Value = 10873.1535
$ws.Range("M136").Value = -75512.181
$ws.Range("N136").Value = -15973.1535

